$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the last existing header cell (G1) to the new
# header cell (H1) so the new "Save" column header matches the other
# header columns (bold font, border, centered alignment).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text and the data value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
